{"js": "// Translate the English subtitle text runs to Swahili.\n// Each (old, new) pair below corresponds to one w:t run text that must be\n// replaced in place, preserving all paragraph/run formatting.\n// Sorted with longest source text first so that a short phrase (e.g.\n// \"velocity\") cannot accidentally match inside a longer sentence that\n// contains it as a substring (e.g. \"As before the ants velocity is one\")\n// before that longer sentence gets its own, full replacement.\nconst replacements = [\n  [\n    \"The dialogue starts at 40 seconds in so I added 27 seconds to the times as they were - John Argentino\",\n    \"Mazungumzo huanza kwa sekunde 40 kwa hivyo niliongeza sekunde 27 kwa nyakati kama zilivyokuwa - John Argentino\"\n  ],\n  [\n    \"with a velocity, let's call it V, which is\",\n    \"kwa kasi, tuiite V, ambayo ni\"\n  ],\n  [\n    \"where exactly to place the two ants on the\",\n    \"wapi hasa kuweka mchwa wawili kwenye\"\n  ],\n  [\n    \"are two ants on a very high stool: a sort\",\n    \"ni mchwa wawili kwenye kinyesi cha juu sana: aina\"\n  ],\n  [\n    \"falling the longest possible. Ants cannot\",\n    \"kuanguka kwa muda mrefu iwezekanavyo. Mchwa hawawezi\"\n  ],\n  [\n    \"probably going to discuss in a different\",\n    \"pengine kwenda kujadili katika tofauti\"\n  ],\n  [\n    \"steep cliffs to both the sides. The flat\",\n    \"miamba mikali kwa pande zote mbili. Gorofa\"\n  ],\n  [\n    \"peak is one meter wide the two ants move\",\n    \"kilele ni mita moja upana wa mchwa wawili hoja\"\n  ],\n  [\n    \"be still: they must move to the right or\",\n    \"tulia: lazima wahamie kulia au\"\n  ],\n  [\n    \"to the left but they must move and after\",\n    \"upande wa kushoto lakini lazima wasogee na baada\"\n  ],\n  [\n    \"order to get the longest time before the\",\n    \"ili kupata muda mrefu zaidi kabla ya\"\n  ],\n  [\n    \"basically the same but now we have three\",\n    \"kimsingi ni sawa lakini sasa tuna tatu\"\n  ],\n  [\n    \"the peak is one meter wide. So, what are\",\n    \"kilele kina upana wa mita moja. Hivyo, ni nini\"\n  ],\n  [\n    \"started! As I said I'm going to discuss\",\n    \"imeanza! Kama nilivyosema nitajadili\"\n  ],\n  [\n    \"equal to one centimeter per second. You\",\n    \"sawa na sentimita moja kwa sekunde. Wewe\"\n  ],\n  [\n    \"top of the mountain. Your purpose is to\",\n    \"juu ya mlima. Kusudi lako ni\"\n  ],\n  [\n    \"make the time the last ant takes before\",\n    \"fanya wakati mchwa wa mwisho huchukua hapo awali\"\n  ],\n  [\n    \"meeting each other they turn around and\",\n    \"wakikutana wanageuka na\"\n  ],\n  [\n    \"so again what are the precise positions\",\n    \"kwa hivyo tena ni nafasi gani sahihi\"\n  ],\n  [\n    \"to get the longest time before the last\",\n    \"kupata muda mrefu zaidi kabla ya mwisho\"\n  ],\n  [\n    \"ant falls down? I hope you enjoyed this\",\n    \"chungu huanguka chini? Natumaini ulifurahia hili\"\n  ],\n  [\n    \"video. Let me just finish writing down\",\n    \"video. Ngoja nimalizie kuandika\"\n  ],\n  [\n    \"the title and, well, I can even draw a\",\n    \"kichwa na, vizuri, naweza hata kuchora a\"\n  ],\n  [\n    \"little ant right here. okay, let's get\",\n    \"mchwa mdogo hapa. sawa, tupate\"\n  ],\n  [\n    \"keep moving with the same but opposite\",\n    \"endelea kusonga na sawa lakini kinyume\"\n  ],\n  [\n    \"centimeter per second, every ant turns\",\n    \"sentimita kwa sekunde, kila mchwa hugeuka\"\n  ],\n  [\n    \"I should place the three ants in order\",\n    \"Ninapaswa kuweka mchwa watatu kwa mpangilio\"\n  ],\n  [\n    \"versions of a more complicated puzzle\",\n    \"matoleo ya fumbo ngumu zaidi\"\n  ],\n  [\n    \"two puzzles in the first puzzle there\",\n    \"mafumbo mawili katika fumbo la kwanza hapo\"\n  ],\n  [\n    \"of Mountain, flat on the top with two\",\n    \"ya Mlima, gorofa juu na mbili\"\n  ],\n  [\n    \"the same for both of them and that is\",\n    \"sawa kwa wote wawili na hiyo ni\"\n  ],\n  [\n    \"can decide the direction towards each\",\n    \"inaweza kuamua mwelekeo kuelekea kila mmoja\"\n  ],\n  [\n    \"ant moves if it is right or left and\",\n    \"mchwa husogea ikiwa ni kulia au kushoto na\"\n  ],\n  [\n    \"where I should place the two ants in\",\n    \"ambapo ninapaswa kuwaweka mchwa wawili ndani\"\n  ],\n  [\n    \"last ant falls? The second puzzle is\",\n    \"chungu mwisho huanguka? Fumbo la pili ni\"\n  ],\n  [\n    \"around after meeting another ant and\",\n    \"karibu baada ya kukutana na mchwa mwingine na\"\n  ],\n  [\n    \"known as the ants puzzle, which I'm\",\n    \"inayojulikana kama fumbo la mchwa, ambalo mimi ni\"\n  ],\n  [\n    \"As before the ants velocity is one\",\n    \"Kama kabla ya mchwa kasi ni moja\"\n  ],\n  [\n    \"okay so the puzzles I'm going to\",\n    \"sawa kwa hivyo mafumbo nitaenda\"\n  ],\n  [\n    \"challenge you with are two basic\",\n    \"changamoto uliyonayo ni mbili za msingi\"\n  ],\n  [\n    \"video do your best and good luck\",\n    \"video fanya bora na bahati nzuri\"\n  ],\n  [\n    \"The ants problem - subtitles:\",\n    \"Tatizo la mchwa - manukuu:\"\n  ],\n  [\n    \"now the precise positions\",\n    \"sasa nafasi sahihi\"\n  ],\n  [\n    \"ants instead of two.\",\n    \"mchwa badala ya wawili.\"\n  ],\n  [\n    \"velocity\",\n    \"kasi\"\n  ],\n  [\n    \"[Music]\",\n    \"[Muziki]\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the English subtitle text runs to Swahili using Word's Find/Replace.\n# Pairs are ordered with the longest source text first so a short phrase\n# (e.g. \"velocity\") cannot accidentally match inside a longer sentence that\n# contains it as a substring (e.g. \"As before the ants velocity is one\")\n# before that longer sentence receives its own, full replacement.\n$pairs = @(\n    @('The dialogue starts at 40 seconds in so I added 27 seconds to the times as they were - John Argentino', 'Mazungumzo huanza kwa sekunde 40 kwa hivyo niliongeza sekunde 27 kwa nyakati kama zilivyokuwa - John Argentino'),\n    @('with a velocity, let''s call it V, which is', 'kwa kasi, tuiite V, ambayo ni'),\n    @('where exactly to place the two ants on the', 'wapi hasa kuweka mchwa wawili kwenye'),\n    @('are two ants on a very high stool: a sort', 'ni mchwa wawili kwenye kinyesi cha juu sana: aina'),\n    @('falling the longest possible. Ants cannot', 'kuanguka kwa muda mrefu iwezekanavyo. Mchwa hawawezi'),\n    @('probably going to discuss in a different', 'pengine kwenda kujadili katika tofauti'),\n    @('steep cliffs to both the sides. The flat', 'miamba mikali kwa pande zote mbili. Gorofa'),\n    @('peak is one meter wide the two ants move', 'kilele ni mita moja upana wa mchwa wawili hoja'),\n    @('be still: they must move to the right or', 'tulia: lazima wahamie kulia au'),\n    @('to the left but they must move and after', 'upande wa kushoto lakini lazima wasogee na baada'),\n    @('order to get the longest time before the', 'ili kupata muda mrefu zaidi kabla ya'),\n    @('basically the same but now we have three', 'kimsingi ni sawa lakini sasa tuna tatu'),\n    @('the peak is one meter wide. So, what are', 'kilele kina upana wa mita moja. Hivyo, ni nini'),\n    @('started! As I said I''m going to discuss', 'imeanza! Kama nilivyosema nitajadili'),\n    @('equal to one centimeter per second. You', 'sawa na sentimita moja kwa sekunde. Wewe'),\n    @('top of the mountain. Your purpose is to', 'juu ya mlima. Kusudi lako ni'),\n    @('make the time the last ant takes before', 'fanya wakati mchwa wa mwisho huchukua hapo awali'),\n    @('meeting each other they turn around and', 'wakikutana wanageuka na'),\n    @('so again what are the precise positions', 'kwa hivyo tena ni nafasi gani sahihi'),\n    @('to get the longest time before the last', 'kupata muda mrefu zaidi kabla ya mwisho'),\n    @('ant falls down? I hope you enjoyed this', 'chungu huanguka chini? Natumaini ulifurahia hili'),\n    @('video. Let me just finish writing down', 'video. Ngoja nimalizie kuandika'),\n    @('the title and, well, I can even draw a', 'kichwa na, vizuri, naweza hata kuchora a'),\n    @('little ant right here. okay, let''s get', 'mchwa mdogo hapa. sawa, tupate'),\n    @('keep moving with the same but opposite', 'endelea kusonga na sawa lakini kinyume'),\n    @('centimeter per second, every ant turns', 'sentimita kwa sekunde, kila mchwa hugeuka'),\n    @('I should place the three ants in order', 'Ninapaswa kuweka mchwa watatu kwa mpangilio'),\n    @('versions of a more complicated puzzle', 'matoleo ya fumbo ngumu zaidi'),\n    @('two puzzles in the first puzzle there', 'mafumbo mawili katika fumbo la kwanza hapo'),\n    @('of Mountain, flat on the top with two', 'ya Mlima, gorofa juu na mbili'),\n    @('the same for both of them and that is', 'sawa kwa wote wawili na hiyo ni'),\n    @('can decide the direction towards each', 'inaweza kuamua mwelekeo kuelekea kila mmoja'),\n    @('ant moves if it is right or left and', 'mchwa husogea ikiwa ni kulia au kushoto na'),\n    @('where I should place the two ants in', 'ambapo ninapaswa kuwaweka mchwa wawili ndani'),\n    @('last ant falls? The second puzzle is', 'chungu mwisho huanguka? Fumbo la pili ni'),\n    @('around after meeting another ant and', 'karibu baada ya kukutana na mchwa mwingine na'),\n    @('known as the ants puzzle, which I''m', 'inayojulikana kama fumbo la mchwa, ambalo mimi ni'),\n    @('As before the ants velocity is one', 'Kama kabla ya mchwa kasi ni moja'),\n    @('okay so the puzzles I''m going to', 'sawa kwa hivyo mafumbo nitaenda'),\n    @('challenge you with are two basic', 'changamoto uliyonayo ni mbili za msingi'),\n    @('video do your best and good luck', 'video fanya bora na bahati nzuri'),\n    @('The ants problem - subtitles:', 'Tatizo la mchwa - manukuu:'),\n    @('now the precise positions', 'sasa nafasi sahihi'),\n    @('ants instead of two.', 'mchwa badala ya wawili.'),\n    @('velocity', 'kasi'),\n    @('[Music]', '[Muziki]')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n"}
